$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.302.06"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "2.522.32"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'535.76"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "'139.48"
$ws.Range("E6").Value = "  -3.69%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E8").Value = "  -2.24%  "

$ws.Range("D9").Value = "2.528.54"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "'5.45"
$ws.Range("E12").Value = "  -2.41%  "

$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "2.968.69"
$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "'23.42"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("D16").Value = "59.224.79"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "2.520.61"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").Value = "'11.06"
$ws.Range("E19").Value = "  -1.96%  "

$ws.Range("D20").Value = "'4.32"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'325.04"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").Value = "'62.90"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").Value = "'0.428"
$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("D28").Value = "'7.84"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("D29").Value = "'6.84"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").Value = "0.0₃0777"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").Value = "'165.15"
$ws.Range("E32").Value = "  +5.14%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("E35").Value = "  -9.25%  "

$ws.Range("D36").Value = "'18.53"
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").Value = "'4.28"
$ws.Range("E37").Value = "  -2.92%  "

$ws.Range("D38").Value = "'1.59"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").Value = "'0.814"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("D42").Value = "'5.29"
$ws.Range("E42").Value = "  -6.73%  "

$ws.Range("D43").Value = "'280.65"
$ws.Range("E43").Value = "  -6.33%  "

$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'10.87"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.598"
$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").Value = "'122.36"
$ws.Range("E48").Value = "  -1.26%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0514"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0225"
$ws.Range("E50").Value = "  -1.94%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'17.78"
$ws.Range("E51").Value = "  -3.36%  "
